# Updated cryptos list on Thu Oct 19 14:47:45 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without Excel coercing
# numeric-looking strings (e.g. "210.34") into a Double. The sheet stores
# every Price/Volume cell as text, so force text format, write the value,
# then drop back to the default "Normal" style so no stray number format
# sticks to the cell (matches the original unstyled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.544.98"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.555.22"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.54%  "
Set-TextValue $ws.Range("D5") "210.34"
$ws.Range("E5").Value = "  -1.02%  "
Set-TextValue $ws.Range("D6") "0.484"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  -0.49%  "
Set-TextValue $ws.Range("D8") "24.32"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("E9").Value = "  -1.15%  "
Set-TextValue $ws.Range("D10") "0.0583"
$ws.Range("E10").Value = "  -0.95%  "
Set-TextValue $ws.Range("D11") "0.0893"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.776.81"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "1.552.98"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "28.526.69"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.510"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "3.62"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("E17").Value = "  -0.88%  "
Set-TextValue $ws.Range("D18") "229.29"
$ws.Range("E18").Value = "  -0.91%  "
Set-TextValue $ws.Range("D19") "7.36"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "0.0₃0671"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  -1.35%  "
Set-TextValue $ws.Range("D23") "8.91"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  +0.73%  "
Set-TextValue $ws.Range("D25") "150.88"
$ws.Range("E25").Value = "  -0.62%  "
Set-TextValue $ws.Range("D26") "14.74"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -0.95%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.48%  "
Set-TextValue $ws.Range("D29") "6.23"
$ws.Range("E29").Value = "  -2.33%  "
Set-TextValue $ws.Range("D30") "0.0459"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "1.390.96"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("E37").Value = "  -2.89%  "
Set-TextValue $ws.Range("D38") "2.65"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  -0.40%  "
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  -0.48%  "
Set-TextValue $ws.Range("D43") "0.775"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("E44").Value = "  +1.34%  "
Set-TextValue $ws.Range("D45") "64.26"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "1.689.45"
$ws.Range("E47").Value = "  -1.47%  "
Set-TextValue $ws.Range("D48") "0.870"
$ws.Range("E48").Value = "  -6.28%  "
Set-TextValue $ws.Range("D49") "43.75"
$ws.Range("E49").Value = "  +5.73%  "
Set-TextValue $ws.Range("D50") "85.25"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  +0.90%  "
